$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 833.48736
$ws.Range("E3").Value = 772.713928
$ws.Range("E4").Value = 797.587807
$ws.Range("E5").Value = 835.976837
$ws.Range("E6").Value = 944.349933594219
$ws.Range("E7").Value = 964.317918498268
$ws.Range("E8").Value = 976.870623402674
$ws.Range("E9").Value = 1004.44940786594
$ws.Range("E10").Value = 1087.68039492592
$ws.Range("E11").Value = 1126.0947600683
$ws.Range("E12").Value = 1149.89261548994
$ws.Range("E13").Value = 1190.00947170954
$ws.Range("E29").Value = 254.1
$ws.Range("E30").Value = 226.104562798092
$ws.Range("E31").Value = 223.630218600954
$ws.Range("E32").Value = 226.205556438792
$ws.Range("E33").Value = 250.266878506616
$ws.Range("E34").Value = 255.502126108265
$ws.Range("E35").Value = 260.307532327884
$ws.Range("E36").Value = 265.344391599365
$ws.Range("E37").Value = 288.933731834977
$ws.Range("E41").Value = 3422.5
$ws.Range("E42").Value = 3475.0001727918
$ws.Range("E43").Value = 3544.60504209738
$ws.Range("E44").Value = 3580.44946756899
$ws.Range("E45").Value = 3665.14221657143
$ws.Range("E46").Value = 3698.55945820078
$ws.Range("E47").Value = 3736.12351960347
$ws.Range("E48").Value = 3772.67679205759
$ws.Range("E49").Value = 3854.46108192878
$ws.Range("E77").Value = 1557.1
$ws.Range("E78").Value = 1611.61604692514
$ws.Range("E79").Value = 1657.01079302431
$ws.Range("E80").Value = 1694.45397666319
$ws.Range("E81").Value = 1685.93358290488
$ws.Range("E82").Value = 1661.02549111116
$ws.Range("E83").Value = 1639.81448236417
$ws.Range("E84").Value = 1642.73245479907
$ws.Range("E85").Value = 1648.78101932228
$ws.Range("E89").Value = 2008.9
$ws.Range("E90").Value = 1725.23333333333
$ws.Range("E91").Value = 1701.4
$ws.Range("E92").Value = 1677.56666666667
$ws.Range("E93").Value = 1732.45226943615
$ws.Range("E94").Value = 1723.61893610282
$ws.Range("E95").Value = 1719.78560276948
$ws.Range("E96").Value = 1715.95226943615
$ws.Range("E97").Value = 1787.82528072634
$ws.Range("E137").Value = 700.1
$ws.Range("E138").Value = 709.317033670381
$ws.Range("E139").Value = 716.556958472619
$ws.Range("E140").Value = 723.601712301784
$ws.Range("E141").Value = 730.722545030267
$ws.Range("E142").Value = 737.906900714206
$ws.Range("E143").Value = 745.392187430077
$ws.Range("E144").Value = 752.35382628783
$ws.Range("E145").Value = 758.700123139062
$ws.Range("E150").Value = 537.804464289892
$ws.Range("E151").Value = 543.293778228426
$ws.Range("E152").Value = 548.635113455558
$ws.Range("E153").Value = 553.435140096838
$ws.Range("E154").Value = 558.876432310261
$ws.Range("E155").Value = 507.751066023249
$ws.Range("E156").Value = 509.847048815395
$ws.Range("E157").Value = 514.147739005903
$ws.Range("E197").Value = 101.2
$ws.Range("E198").Value = 90.0503020667727
$ws.Range("E199").Value = 89.0648489666137
$ws.Range("E200").Value = 90.0905246422893
$ws.Range("E201").Value = 91.3334022257552
$ws.Range("E202").Value = 93.2439745627981
$ws.Range("E203").Value = 94.9976788553259
$ws.Range("E204").Value = 96.8358505564388
$ws.Range("E205").Value = 98.1873290937997
$ws.Range("E209").Value = 1951.9
$ws.Range("E210").Value = 1978.53353532519
$ws.Range("E211").Value = 2011.95360806146
$ws.Range("E212").Value = 2035.53278727767
$ws.Range("E213").Value = 2061.48596715345
$ws.Range("E214").Value = 2083.69984430965
$ws.Range("E215").Value = 2106.69406718421
$ws.Range("E216").Value = 2129.50119350236
$ws.Range("E217").Value = 2151.20880788586
$ws.Range("E221").Value = 2389.6
$ws.Range("E222").Value = 2441.64114099415
$ws.Range("E223").Value = 2506.18184029607
$ws.Range("E224").Value = 2578.20285177688
$ws.Range("E225").Value = 2634.35141978861
$ws.Range("E226").Value = 2685.52777121801
$ws.Range("E227").Value = 2731.33891384639
$ws.Range("E228").Value = 2768.14418216863
$ws.Range("E229").Value = 2802.11549650558
$ws.Range("E234").Value = 103.107612369063
$ws.Range("E235").Value = 90.6057100111654
$ws.Range("E236").Value = 79.619675967694
$ws.Range("E237").Value = 69.9657096712712
$ws.Range("E238").Value = 61.4822965593436
$ws.Range("E239").Value = 54.0275058735411
$ws.Range("E240").Value = 47.476616103598
$ws.Range("E241").Value = 41.7200283486025
